# Fruta / hortaliza, semanal
# Insert two new weekly price rows (date 44448) for
# "Terminal La Palmera de La Serena" / Cebolla, pushing the existing
# rows 371-379 down to 373-381.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 371; Excel shifts the
# existing rows 371-379 down to 373-381 and copies formatting (incl.
# the date-number-format style on column D) from the row above.
$ws.Rows("371:372").Insert()

# New row 371: "1a (guarda)"
$ws.Cells.Item(371, 1).Value = 8
$ws.Cells.Item(371, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(371, 3).Value = "Coquimbo"
$ws.Cells.Item(371, 4).Value = 44448
$ws.Cells.Item(371, 5).Value = 4
$ws.Cells.Item(371, 6).Value = 100112004
$ws.Cells.Item(371, 7).Value = "Cebolla"
$ws.Cells.Item(371, 8).Value = "Sin especificar"
$ws.Cells.Item(371, 9).Value = "1a (guarda)"
$ws.Cells.Item(371, 10).Value = 2560
$ws.Cells.Item(371, 11).Value = 5300
$ws.Cells.Item(371, 12).Value = 5500
$ws.Cells.Item(371, 13).Value = 5400
$ws.Cells.Item(371, 14).Value = "$/malla 16 kilos"
$ws.Cells.Item(371, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(371, 16).Value = 338
$ws.Cells.Item(371, 17).Value = 16
$ws.Cells.Item(371, 18).Value = "Hortaliza"

# New row 372: "2a (guarda)"
$ws.Cells.Item(372, 1).Value = 8
$ws.Cells.Item(372, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(372, 3).Value = "Coquimbo"
$ws.Cells.Item(372, 4).Value = 44448
$ws.Cells.Item(372, 5).Value = 4
$ws.Cells.Item(372, 6).Value = 100112004
$ws.Cells.Item(372, 7).Value = "Cebolla"
$ws.Cells.Item(372, 8).Value = "Sin especificar"
$ws.Cells.Item(372, 9).Value = "2a (guarda)"
$ws.Cells.Item(372, 10).Value = 1500
$ws.Cells.Item(372, 11).Value = 4800
$ws.Cells.Item(372, 12).Value = 5000
$ws.Cells.Item(372, 13).Value = 4900
$ws.Cells.Item(372, 14).Value = "$/malla 16 kilos"
$ws.Cells.Item(372, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(372, 16).Value = 306
$ws.Cells.Item(372, 17).Value = 16
$ws.Cells.Item(372, 18).Value = "Hortaliza"
